# Fix spi flash layout size
# Update part manufacturer numbers, resistor value, and KiCost scrape timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manufacturer part numbers (column D) for the capacitor rows
$ws.Range("D7").Value = "C1608X5R1A106K080AC"
$ws.Range("D8").Value = "C0402C104J8RACTU"
$ws.Range("D9").Value = "C0402C103K8RACAUTO"
$ws.Range("D10").Value = "CC0402KRX5R6BB105"

# Resistor value for R2,R3 changed from 56 to 68.
# This column stores the "Value" as text (shared string), so force a
# text entry (leading apostrophe, like typing '68 into the cell) rather
# than letting it be auto-recognized as a number.
$ws.Range("B17").Value = "'68"

# KiCost scrape date/time metadata
$ws.Range("B3").Value = "wo 07 nov 2018 08:29:51 CET"
$ws.Range("B4").Value = "2018-11-07 08:31:08"
